$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '304.14'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '1.77%'

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '31.71'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '0.04%'

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.185'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '1.27%'

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07477'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-0.52%'

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '2.433'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '41.39%'

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '8.020'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '3.06%'

# Row 8
$ws.Range("B8").Value = 'GateToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.871'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '1.94%'

# Row 9
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9149'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-1.12%'

# Row 10
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1733'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '1.33%'

# Row 11
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07686'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '4.20%'

# Row 12
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08178'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '3.00%'

# Row 13
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03042'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '0.15%'

# Row 14
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09958'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.56%'

# Row 15
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001515'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '1.51%'

# Row 16
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006140'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-0.41%'

# Row 17
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.500'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '1.60%'

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.236'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '0.53%'

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3261'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-0.93%'

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1340'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '0.36%'

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.660'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '2.08%'

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.04603'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-1.19%'

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.1567'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '1.11%'

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001264'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '3.99%'

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004527'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '2.35%'

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001300'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-7.16%'

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0002744'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '51.73%'

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01757'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '5.28%'

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04530'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-0.45%'

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007396'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '4.29%'

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1364'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '2.73%'

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002180'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '5.80%'

# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-14.69%'

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00006490'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '6.59%'
